# "Make Character UI." -- adds new UI localization strings (upgrade /
# def-shred / skills / empty-relic) to the "UI" sheet, and moves the
# active-tab/selection focus from the "STR" sheet back to "UI".

$wb = $excel.ActiveWorkbook

# --- 1. Append the new localization rows to the "UI" sheet (sheet1) -----
$wsUI = $wb.Worksheets.Item(1)

$newRows = @(
    @("UI_UPGRADE_LEVEL",     "Upgrade",                                    "Tăng cấp"),
    @("UI_UPGRADE_MAX_LEVEL", "Upgrade to Lv.10",                           "Tăng tới cấp 10"),
    @("UI_DEF_SHRED",         "DEF Shred",                                  "Khả năng giảm giáp"),
    @("UI_SKILLS",            "Skill",                                      "Kĩ năng"),
    @("UI_WEAPON_EMPTY",      "No Relic equipped. Please go to change it.", "Chưa trang bị Pháp Bảo, vui lòng đến để thay đổi.")
)

$startRow = 56
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $wsUI.Cells.Item($r, 1).Value = $row[0]
    $wsUI.Cells.Item($r, 2).Value = $row[1]
    $wsUI.Cells.Item($r, 3).Value = $row[2]
}

# --- 2. Move the STR sheet's selection (it no longer stays the active tab)
$wsSTR = $wb.Worksheets.Item(3)
$wsSTR.Select()
$wsSTR.Range("A149").Select()

# --- 3. Re-activate the UI sheet and select the newly-added last row ----
$wsUI.Activate()
$wsUI.Range("C60").Select()

Write-Host "Character UI localization rows added."
